$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new shared string + data rows 197-213 (Keel Ridge Mine) ---
$ws.Cells.Item(197,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(197,4).Value = 1996
$ws.Cells.Item(197,6).Value = 242
$ws.Cells.Item(197,7).Formula = "=(F197-I197)/(J197-I197)"
$ws.Cells.Item(197,8).Value = 1996
$ws.Cells.Item(197,9).Value = 18
$ws.Cells.Item(197,10).Value = 298.8333333333333
$ws.Cells.Item(197,11).Value = 7.97569444444444
$ws.Cells.Item(197,12).Value = 2024
$ws.Cells.Item(197,13).Value = 7.22222222222222
$ws.Cells.Item(197,14).Value = 35069

$ws.Cells.Item(198,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(198,4).Value = 2004
$ws.Cells.Item(198,6).Value = 393
$ws.Cells.Item(198,7).Formula = "=(F198-I198)/(J198-I198)"
$ws.Cells.Item(198,8).Value = 2005
$ws.Cells.Item(198,9).Value = 18
$ws.Cells.Item(198,10).Value = 298.8333333333333
$ws.Cells.Item(198,11).Value = 7.97569444444444
$ws.Cells.Item(198,12).Value = 2024
$ws.Cells.Item(198,13).Value = 8.88888888888889
$ws.Cells.Item(198,14).Value = 38350

$ws.Cells.Item(199,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(199,4).Value = 2007
$ws.Cells.Item(199,6).Value = 267
$ws.Cells.Item(199,7).Formula = "=(F199-I199)/(J199-I199)"
$ws.Cells.Item(199,8).Value = 2007
$ws.Cells.Item(199,9).Value = 18
$ws.Cells.Item(199,10).Value = 298.833333333333
$ws.Cells.Item(199,11).Value = 7.97569444444444
$ws.Cells.Item(199,12).Value = 2024
$ws.Cells.Item(199,13).Value = 7.5
$ws.Cells.Item(199,14).Value = 39165

$ws.Cells.Item(200,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(200,4).Value = 2009
$ws.Cells.Item(200,6).Value = 203
$ws.Cells.Item(200,7).Formula = "=(F200-I200)/(J200-I200)"
$ws.Cells.Item(200,8).Value = 2009
$ws.Cells.Item(200,9).Value = 18
$ws.Cells.Item(200,10).Value = 298.833333333333
$ws.Cells.Item(200,11).Value = 7.97569444444444
$ws.Cells.Item(200,12).Value = 2024
$ws.Cells.Item(200,13).Value = 7.97569444444444
$ws.Cells.Item(200,14).Value = 39893

$ws.Cells.Item(201,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(201,4).Value = 2009
$ws.Cells.Item(201,6).Value = 275
$ws.Cells.Item(201,7).Formula = "=(F201-I201)/(J201-I201)"
$ws.Cells.Item(201,8).Value = 2010
$ws.Cells.Item(201,9).Value = 18
$ws.Cells.Item(201,10).Value = 298.833333333333
$ws.Cells.Item(201,11).Value = 7.97569444444444
$ws.Cells.Item(201,12).Value = 2024
$ws.Cells.Item(201,13).Value = 7.97569444444444
$ws.Cells.Item(201,14).Value = 40131

$ws.Cells.Item(202,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(202,4).Value = 2010
$ws.Cells.Item(202,6).Value = 290
$ws.Cells.Item(202,7).Formula = "=(F202-I202)/(J202-I202)"
$ws.Cells.Item(202,8).Value = 2010
$ws.Cells.Item(202,9).Value = 18
$ws.Cells.Item(202,10).Value = 298.833333333333
$ws.Cells.Item(202,11).Value = 7.97569444444444
$ws.Cells.Item(202,12).Value = 2024
$ws.Cells.Item(202,13).Value = 7.97569444444444
$ws.Cells.Item(202,14).Value = 40264

$ws.Cells.Item(203,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(203,4).Value = 2010
$ws.Cells.Item(203,6).Value = 324
$ws.Cells.Item(203,7).Formula = "=(F203-I203)/(J203-I203)"
$ws.Cells.Item(203,8).Value = 2011
$ws.Cells.Item(203,9).Value = 18
$ws.Cells.Item(203,10).Value = 298.833333333333
$ws.Cells.Item(203,11).Value = 7.97569444444444
$ws.Cells.Item(203,12).Value = 2024
$ws.Cells.Item(203,13).Value = 9
$ws.Cells.Item(203,14).Value = 40488

$ws.Cells.Item(204,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(204,4).Value = 2011
$ws.Cells.Item(204,6).Value = 262
$ws.Cells.Item(204,7).Formula = "=(F204-I204)/(J204-I204)"
$ws.Cells.Item(204,8).Value = 2011
$ws.Cells.Item(204,9).Value = 18
$ws.Cells.Item(204,10).Value = 298.833333333333
$ws.Cells.Item(204,11).Value = 7.97569444444444
$ws.Cells.Item(204,12).Value = 2024
$ws.Cells.Item(204,13).Value = 7.97569444444444
$ws.Cells.Item(204,14).Value = 40628

$ws.Cells.Item(205,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(205,4).Value = 2013
$ws.Cells.Item(205,6).Value = 222
$ws.Cells.Item(205,7).Formula = "=(F205-I205)/(J205-I205)"
$ws.Cells.Item(205,8).Value = 2013
$ws.Cells.Item(205,9).Value = 18
$ws.Cells.Item(205,10).Value = 298.833333333333
$ws.Cells.Item(205,11).Value = 7.97569444444444
$ws.Cells.Item(205,12).Value = 2024
$ws.Cells.Item(205,13).Value = 8.63888888888889
$ws.Cells.Item(205,14).Value = 41348

$ws.Cells.Item(206,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(206,4).Value = 2014
$ws.Cells.Item(206,6).Value = 300
$ws.Cells.Item(206,7).Formula = "=(F206-I206)/(J206-I206)"
$ws.Cells.Item(206,8).Value = 2014
$ws.Cells.Item(206,9).Value = 18
$ws.Cells.Item(206,10).Value = 298.833333333333
$ws.Cells.Item(206,11).Value = 7.97569444444444
$ws.Cells.Item(206,12).Value = 2024
$ws.Cells.Item(206,13).Value = 7.97569444444444
$ws.Cells.Item(206,14).Value = 41694

$ws.Cells.Item(207,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(207,4).Value = 2015
$ws.Cells.Item(207,6).Value = 587
$ws.Cells.Item(207,7).Formula = "=(F207-I207)/(J207-I207)"
$ws.Cells.Item(207,8).Value = 2015
$ws.Cells.Item(207,9).Value = 18
$ws.Cells.Item(207,10).Value = 298.833333333333
$ws.Cells.Item(207,11).Value = 7.97569444444444
$ws.Cells.Item(207,12).Value = 2024
$ws.Cells.Item(207,13).Value = 7.97569444444444
$ws.Cells.Item(207,14).Value = 42057

$ws.Cells.Item(208,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(208,4).Value = 2016
$ws.Cells.Item(208,6).Value = 221
$ws.Cells.Item(208,7).Formula = "=(F208-I208)/(J208-I208)"
$ws.Cells.Item(208,8).Value = 2016
$ws.Cells.Item(208,9).Value = 18
$ws.Cells.Item(208,10).Value = 298.833333333333
$ws.Cells.Item(208,11).Value = 7.97569444444444
$ws.Cells.Item(208,12).Value = 2024
$ws.Cells.Item(208,13).Value = 8.05555555555556
$ws.Cells.Item(208,14).Value = 42420

$ws.Cells.Item(209,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(209,4).Value = 2017
$ws.Cells.Item(209,5).Value = 0
$ws.Cells.Item(209,6).Value = 45
$ws.Cells.Item(209,7).Formula = "=(F209-I209)/(J209-I209)"
$ws.Cells.Item(209,8).Value = 2017
$ws.Cells.Item(209,9).Value = 18
$ws.Cells.Item(209,10).Value = 298.833333333333
$ws.Cells.Item(209,11).Value = 7.97569444444444
$ws.Cells.Item(209,12).Value = 2024
$ws.Cells.Item(209,13).Value = 8.25
$ws.Cells.Item(209,14).Value = 42777

$ws.Cells.Item(210,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(210,4).Value = 2019
$ws.Cells.Item(210,5).Value = 2
$ws.Cells.Item(210,6).Value = 103
$ws.Cells.Item(210,7).Formula = "=(F210-I210)/(J210-I210)"
$ws.Cells.Item(210,8).Value = 2019
$ws.Cells.Item(210,9).Value = 18
$ws.Cells.Item(210,10).Value = 298.833333333333
$ws.Cells.Item(210,11).Value = 7.97569444444444
$ws.Cells.Item(210,12).Value = 2024
$ws.Cells.Item(210,13).Value = 7.97569444444444
$ws.Cells.Item(210,14).Value = 43479

$ws.Cells.Item(211,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(211,4).Value = 2021
$ws.Cells.Item(211,5).Value = 4
$ws.Cells.Item(211,6).Value = 107
$ws.Cells.Item(211,7).Formula = "=(F211-I211)/(J211-I211)"
$ws.Cells.Item(211,8).Value = 2021
$ws.Cells.Item(211,9).Value = 18
$ws.Cells.Item(211,10).Value = 298.833333333333
$ws.Cells.Item(211,11).Value = 7.97569444444444
$ws.Cells.Item(211,12).Value = 2024
$ws.Cells.Item(211,13).Value = 7.97569444444444
$ws.Cells.Item(211,14).Value = 44258

$ws.Cells.Item(212,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(212,4).Value = 2023
$ws.Cells.Item(212,5).Value = 6
$ws.Cells.Item(212,6).Value = 30
$ws.Cells.Item(212,7).Formula = "=(F212-I212)/(J212-I212)"
$ws.Cells.Item(212,8).Value = 2023
$ws.Cells.Item(212,9).Value = 18
$ws.Cells.Item(212,10).Value = 298.833333333333
$ws.Cells.Item(212,11).Value = 7.97569444444444
$ws.Cells.Item(212,12).Value = 2024
$ws.Cells.Item(212,13).Value = 7.97569444444444
$ws.Cells.Item(212,14).Value = 45028

$ws.Cells.Item(213,1).Value = "Keel Ridge Mine"
$ws.Cells.Item(213,4).Value = 2024
$ws.Cells.Item(213,5).Value = 7
$ws.Cells.Item(213,6).Value = 18
$ws.Cells.Item(213,7).Formula = "=(F213-I213)/(J213-I213)"
$ws.Cells.Item(213,8).Value = 2024
$ws.Cells.Item(213,9).Value = 18
$ws.Cells.Item(213,10).Value = 298.833333333333
$ws.Cells.Item(213,11).Value = 7.97569444444444
$ws.Cells.Item(213,12).Value = 2024
$ws.Cells.Item(213,13).Value = 6.25
$ws.Cells.Item(213,14).Value = 45394

# --- Apply the date number format to column N (reuse same style via copy/paste-format) ---
$ws.Cells.Item(197,14).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(197,14).Copy()
$ws.Range("N198:N213").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column N width ---
$ws.Columns.Item(14).ColumnWidth = 10.83

# --- Sheet view: selection on A214 (new first empty row) ---
$ws.Range("A214").Select()
